$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 139
$ws.Range("I53").Value = 102.166664
$ws.Range("K53").Value = 102.166664
$ws.Range("M53").Value = 534.833336

$ws.Range("H62").Value = 3263.182
$ws.Range("I62").Value = 2985.375
$ws.Range("J62").Value = 4004
$ws.Range("K62").Value = 2985.375
$ws.Range("L62").Value = 4004
$ws.Range("M62").Value = -2361.375
$ws.Range("N62").Value = -5252

$ws.Range("H65").Value = 3263.182
$ws.Range("I65").Value = 2985.375
$ws.Range("J65").Value = 4004
$ws.Range("K65").Value = 14926.875
$ws.Range("L65").Value = 20020
$ws.Range("M65").Value = -11806.875
$ws.Range("N65").Value = -26260

$ws.Range("H70").Value = 925.5
$ws.Range("I70").Value = 851
$ws.Range("K70").Value = 2553
$ws.Range("M70").Value = -2283

$ws.Range("H73").Value = 925.5
$ws.Range("I73").Value = 851
$ws.Range("K73").Value = 2553
$ws.Range("M73").Value = -1617

$ws.Range("H86").Value = 166668750
$ws.Range("I86").Value = 166668750
$ws.Range("K86").Value = 166668750
$ws.Range("M86").Value = -166667627

$ws.Range("H89").Value = 166668750
$ws.Range("I89").Value = 166668750
$ws.Range("K89").Value = 833343750
$ws.Range("M89").Value = -833338134

$ws.Range("H135").Value = 3120.6667
$ws.Range("I135").Value = 3160.75
$ws.Range("J135").Value = 2800
$ws.Range("K135").Value = 28446.75
$ws.Range("L135").Value = 25200
$ws.Range("M135").Value = -25911.75
$ws.Range("N135").Value = -30270

$ws.Range("H138").Value = 3116.6667
$ws.Range("I138").Value = 3869.3684
$ws.Range("J138").Value = 2669.75
$ws.Range("K138").Value = 11608.1052
$ws.Range("L138").Value = 8009.25
$ws.Range("M138").Value = -6468.1052
$ws.Range("N138").Value = -18289.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 14751.5
$ws.Range("I14").Value = 14751.5
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 14751.5
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -14576.5
$ws.Range("N14").ClearContents()

$ws.Range("H45").Value = 37546.5
$ws.Range("I45").Value = 54002
$ws.Range("K45").Value = 54002
$ws.Range("M45").Value = -53625

$ws.Range("H61").Value = 3783.3333
$ws.Range("J61").Value = 6666.6665
$ws.Range("L61").Value = 6666.6665
$ws.Range("N61").Value = -7090.6665

$ws.Range("H63").Value = 4699.875
$ws.Range("I63").Value = 3749.75
$ws.Range("J63").Value = 5650
$ws.Range("K63").Value = 3749.75
$ws.Range("L63").Value = 5650
$ws.Range("M63").Value = -3063.75
$ws.Range("N63").Value = -7022

$ws.Range("H66").Value = 4699.875
$ws.Range("I66").Value = 3749.75
$ws.Range("J66").Value = 5650
$ws.Range("K66").Value = 18748.75
$ws.Range("L66").Value = 28250
$ws.Range("M66").Value = -15316.75
$ws.Range("N66").Value = -35114

$ws.Range("H74").Value = 294534.6
$ws.Range("I74").Value = 556644.3
$ws.Range("K74").Value = 556644.3
$ws.Range("M74").Value = -555770.3

$ws.Range("H77").Value = 294534.6
$ws.Range("I77").Value = 556644.3
$ws.Range("K77").Value = 2783221.5
$ws.Range("M77").Value = -2778853.5

$ws.Range("H122").Value = 1479.1052
$ws.Range("I122").Value = 1479.1052
$ws.Range("K122").Value = 4437.3156
$ws.Range("M122").Value = -1987.3156

$ws.Range("H132").Value = 1087.7407
$ws.Range("I132").Value = 723.7083
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 2171.1249
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = 358.8751000000002
$ws.Range("N132").Value = -17060

$ws.Range("H135").Value = 112330.664
$ws.Range("J135").Value = 112330.664
$ws.Range("L135").Value = 112330.664
$ws.Range("N135").Value = -122470.664

$ws.Range("H136").Value = 3783.3333
$ws.Range("J136").Value = 6666.6665
$ws.Range("L136").Value = 19999.9995
$ws.Range("N136").Value = -25099.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3656.775
$ws.Range("I134").Value = 3447.606
$ws.Range("K134").Value = 10342.818
$ws.Range("M134").Value = -7807.818000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 5134.579
$ws.Range("I122").Value = 4600.4287
$ws.Range("K122").Value = 13801.2861
$ws.Range("M122").Value = -11351.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 2945.6365
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 2945.6365
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 8836.9095
$ws.Range("N21").Value = -9182.9095
$ws.Range("M21").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()

$ws.Range("H80").Value = 76925270
$ws.Range("I80").Value = 100001820
$ws.Range("J80").Value = 3429.3333
$ws.Range("K80").Value = 100001820
$ws.Range("L80").Value = 3429.3333
$ws.Range("M80").Value = -100000822
$ws.Range("N80").Value = -5425.3333

$ws.Range("H83").Value = 76925270
$ws.Range("I83").Value = 100001820
$ws.Range("J83").Value = 3429.3333
$ws.Range("K83").Value = 500009100
$ws.Range("L83").Value = 17146.6665
$ws.Range("M83").Value = -500004108
$ws.Range("N83").Value = -27130.6665

$ws.Range("H113").Value = 3348
$ws.Range("I113").Value = 2668.8125
$ws.Range("K113").Value = 2668.8125
$ws.Range("M113").Value = -498.8125

$ws.Range("H122").Value = 4505.3105
$ws.Range("I122").Value = 3042.3845
$ws.Range("J122").Value = 5693.9375
$ws.Range("K122").Value = 9127.1535
$ws.Range("L122").Value = 17081.8125
$ws.Range("M122").Value = -6677.1535
$ws.Range("N122").Value = -21981.8125

$ws.Range("H126").Value = 2234.6667
$ws.Range("I126").Value = 2234.6667
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6704.000100000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4234.000100000001
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 10000
$ws.Range("J17").Value = 10000
$ws.Range("L17").Value = 10000
$ws.Range("N17").Value = -10340

$ws.Range("H40").Value = 4469.5093
$ws.Range("I40").Value = 4711.36
$ws.Range("K40").Value = 4711.36
$ws.Range("M40").Value = -4575.36

$ws.Range("H122").Value = 3456.4285
$ws.Range("I122").Value = 3599.2307
$ws.Range("K122").Value = 10797.6921
$ws.Range("M122").Value = -8347.6921

$ws.Range("H132").Value = 3596
$ws.Range("I132").Value = 1945.5
$ws.Range("K132").Value = 5836.5
$ws.Range("M132").Value = -3306.5

$ws.Range("H136").Value = 4158.05
$ws.Range("I136").Value = 4436
$ws.Range("K136").Value = 13308
$ws.Range("M136").Value = -10758

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 834127.3
$ws.Range("I2").Value = 834127.3
$ws.Range("K2").Value = 834127.3
$ws.Range("M2").Value = -834015.3

$ws.Range("H100").Value = 90910210
$ws.Range("J100").Value = 250001380
$ws.Range("L100").Value = 500002760
$ws.Range("N100").Value = -500003842

$ws.Range("H113").Value = 650.75
$ws.Range("I113").Value = 437.5
$ws.Range("K113").Value = 1312.5
$ws.Range("M113").Value = 857.5

$ws.Range("H122").Value = 10001455
$ws.Range("I122").Value = 1243.6842
$ws.Range("K122").Value = 3731.0526
$ws.Range("M122").Value = -1281.0526

$ws.Range("H126").Value = 2263.5833
$ws.Range("I126").Value = 1394.7142
$ws.Range("K126").Value = 4184.142599999999
$ws.Range("M126").Value = -1714.142599999999

$ws.Range("H132").Value = 4435.5835
$ws.Range("I132").Value = 7160.6
$ws.Range("K132").Value = 21481.8
$ws.Range("M132").Value = -18951.8

$ws.Range("H139").Value = 80832.336
$ws.Range("J139").Value = 79998.91
$ws.Range("K139").Value = 79998.91
$ws.Range("N139").Value = -90278.91
